$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ucn2"
$ws.Range("C2").Value = "Crhr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1458113333333333
$ws.Range("H2").Value = 0.437434
$ws.Range("I2").Value = 0.1510638746080812
$ws.Range("J2").Value = 0.1510638746080812
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04713733333333334
$ws.Range("N2").Value = 0.141412
$ws.Range("O2").Value = 0.00822544460891999
$ws.Range("P2").Value = 0.00822544460891999
$ws.Range("Q2").Value = 0.006873157423111113
$ws.Range("R2").Value = 0.06185841680800001
$ws.Range("S2").Value = 0.001242567532997607
$ws.Range("T2").Value = 0.001242567532997607

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ucn2"
$ws.Range("C3").Value = "Crhr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1458113333333333
$ws.Range("H3").Value = 0.437434
$ws.Range("I3").Value = 0.1510638746080812
$ws.Range("J3").Value = 0.1510638746080812
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1124846666666667
$ws.Range("N3").Value = 0.337454
$ws.Range("O3").Value = 0.01962852646917154
$ws.Range("P3").Value = 0.01962852646917154
$ws.Range("Q3").Value = 0.01640153922622222
$ws.Range("R3").Value = 0.147613853036
$ws.Range("S3").Value = 0.002965161261280332
$ws.Range("T3").Value = 0.002965161261280332

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ucn2"
$ws.Range("C4").Value = "Crhr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1458113333333333
$ws.Range("H4").Value = 0.437434
$ws.Range("I4").Value = 0.1510638746080812
$ws.Range("J4").Value = 0.1510638746080812
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.571051
$ws.Range("N4").Value = 16.713153
$ws.Range("O4").Value = 0.9721460289219085
$ws.Range("P4").Value = 0.9721460289219085
$ws.Range("Q4").Value = 0.812322374378
$ws.Range("R4").Value = 7.310901369402
$ws.Range("S4").Value = 0.1468561458138033
$ws.Range("T4").Value = 0.1468561458138033

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ucn2"
$ws.Range("C5").Value = "Crhr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5810133333333333
$ws.Range("H5").Value = 1.74304
$ws.Range("I5").Value = 0.6019430954083812
$ws.Range("J5").Value = 0.6019430954083812
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04713733333333334
$ws.Range("N5").Value = 0.141412
$ws.Range("O5").Value = 0.00822544460891999
$ws.Range("P5").Value = 0.00822544460891999
$ws.Range("Q5").Value = 0.02738741916444444
$ws.Range("R5").Value = 0.24648677248
$ws.Range("S5").Value = 0.00495124958900348
$ws.Range("T5").Value = 0.00495124958900348

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ucn2"
$ws.Range("C6").Value = "Crhr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5810133333333333
$ws.Range("H6").Value = 1.74304
$ws.Range("I6").Value = 0.6019430954083812
$ws.Range("J6").Value = 0.6019430954083812
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1124846666666667
$ws.Range("N6").Value = 0.337454
$ws.Range("O6").Value = 0.01962852646917154
$ws.Range("P6").Value = 0.01962852646917154
$ws.Range("Q6").Value = 0.06535509112888888
$ws.Range("R6").Value = 0.5881958201599999
$ws.Range("S6").Value = 0.01181525598115846
$ws.Range("T6").Value = 0.01181525598115846

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ucn2"
$ws.Range("C7").Value = "Crhr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5810133333333333
$ws.Range("H7").Value = 1.74304
$ws.Range("I7").Value = 0.6019430954083812
$ws.Range("J7").Value = 0.6019430954083812
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.571051
$ws.Range("N7").Value = 16.713153
$ws.Range("O7").Value = 0.9721460289219085
$ws.Range("P7").Value = 0.9721460289219085
$ws.Range("Q7").Value = 3.23685491168
$ws.Range("R7").Value = 29.13169420512
$ws.Range("S7").Value = 0.5851765898382193
$ws.Range("T7").Value = 0.5851765898382193

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ucn2"
$ws.Range("C8").Value = "Crhr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.238405
$ws.Range("H8").Value = 0.7152149999999999
$ws.Range("I8").Value = 0.2469930299835376
$ws.Range("J8").Value = 0.2469930299835376
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04713733333333334
$ws.Range("N8").Value = 0.141412
$ws.Range("O8").Value = 0.00822544460891999
$ws.Range("P8").Value = 0.00822544460891999
$ws.Range("Q8").Value = 0.01123777595333333
$ws.Range("R8").Value = 0.10113998358
$ws.Range("S8").Value = 0.002031627486918903
$ws.Range("T8").Value = 0.002031627486918903

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ucn2"
$ws.Range("C9").Value = "Crhr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.238405
$ws.Range("H9").Value = 0.7152149999999999
$ws.Range("I9").Value = 0.2469930299835376
$ws.Range("J9").Value = 0.2469930299835376
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1124846666666667
$ws.Range("N9").Value = 0.337454
$ws.Range("O9").Value = 0.01962852646917154
$ws.Range("P9").Value = 0.01962852646917154
$ws.Range("Q9").Value = 0.02681690695666666
$ws.Range("R9").Value = 0.2413521626099999
$ws.Range("S9").Value = 0.004848109226732747
$ws.Range("T9").Value = 0.004848109226732747

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ucn2"
$ws.Range("C10").Value = "Crhr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.238405
$ws.Range("H10").Value = 0.7152149999999999
$ws.Range("I10").Value = 0.2469930299835376
$ws.Range("J10").Value = 0.2469930299835376
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.571051
$ws.Range("N10").Value = 16.713153
$ws.Range("O10").Value = 0.9721460289219085
$ws.Range("P10").Value = 0.9721460289219085
$ws.Range("Q10").Value = 1.328166413655
$ws.Range("R10").Value = 11.953497722895
$ws.Range("S10").Value = 0.240113293269886
$ws.Range("T10").Value = 0.240113293269886

